# "corrected the expense rules"
# Rum quantity (Actual!B3) 2 -> 4, and cigarette quantity (Actual!B4) 51 -> 56.
# Downstream formulas (D3, D4, E2, H2 on "Actual") recalc automatically.
# The static snapshot copies on "actual_cost_v1" and "current_total_expense_v1"
# are updated to match, and their selections are adjusted as recorded by Excel.

$wb = $excel.ActiveWorkbook

# --- Actual sheet: correct the source quantities -------------------------
$wsActual = $wb.Worksheets.Item("Actual")
$wsActual.Range("B3").Value = 4
$wsActual.Range("B4").Value = 56
$wsActual.Range("I2").Value = 24

# --- actual_cost_v1: static snapshot of the Actual sheet ------------------
$wsActualCost = $wb.Worksheets.Item("actual_cost_v1")
$wsActualCost.Range("B3").Value = 4
$wsActualCost.Range("D3").Value = 2720
$wsActualCost.Range("B4").Value = 56
$wsActualCost.Range("D4").Value = 560
$wsActualCost.Range("E2").Value = 8774
$wsActualCost.Activate()
$wsActualCost.Range("A1:E18").Select()

# --- current_total_expense_v1: static snapshot of totals ------------------
$wsTotals = $wb.Worksheets.Item("current_total_expense_v1")
$wsTotals.Range("B2").Value = 11226
$wsTotals.Range("C2").Value = 24
$wsTotals.Activate()
$wsTotals.Range("A2:C2").Select()
